$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Fill in data row first so shared-string indices land in the expected order
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"

# Fill in header row
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Passsword"

# Auto-fit column B so its width reflects the "pointofsale" content
$ws.Columns("B").AutoFit() | Out-Null

# Select B2 to match final selection state
$ws.Range("B2").Select()
